$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for each data row (2-52).
# Every row's date is being advanced by one day (45177 -> 45178).
for ($r = 2; $r -le 52; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
